$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for first file row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-17 15:07:50"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime" for first file row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-17 15:07:44"
$wsZhCn.Range("K2").Value = "2016-08-17 15:08:06"

# de-de sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime" for first file row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-17 15:07:50"
$wsDeDe.Range("K2").Value = "2016-08-17 15:08:17"
